{"js": "// Update the division-fact answers in the single table on the page.\n// Every populated cell (5 per \"fact row\", at table rows 0, 4, 8, 12, 16)\n// gets its text replaced with a newly generated problem/answer string.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of (tableRowIndex, columnIndex) -> new cell text, in document order.\nconst updates = [\n  [0, 0, \"89\u00f75=17, 4\"],\n  [0, 1, \"33\u00f74=8, 1\"],\n  [0, 2, \"96\u00f78=12, 0\"],\n  [0, 3, \"74\u00f73=24, 2\"],\n  [0, 4, \"50\u00f78=6, 2\"],\n\n  [4, 0, \"89\u00f74=22, 1\"],\n  [4, 1, \"44\u00f74=11, 0\"],\n  [4, 2, \"22\u00f73=7, 1\"],\n  [4, 3, \"31\u00f76=5, 1\"],\n  [4, 4, \"60\u00f79=6, 6\"],\n\n  [8, 0, \"10\u00f73=3, 1\"],\n  [8, 1, \"53\u00f78=6, 5\"],\n  [8, 2, \"24\u00f79=2, 6\"],\n  [8, 3, \"17\u00f73=5, 2\"],\n  [8, 4, \"92\u00f74=23, 0\"],\n\n  [12, 0, \"73\u00f76=12, 1\"],\n  [12, 1, \"94\u00f76=15, 4\"],\n  [12, 2, \"88\u00f78=11, 0\"],\n  [12, 3, \"32\u00f73=10, 2\"],\n  [12, 4, \"92\u00f78=11, 4\"],\n\n  [16, 0, \"93\u00f73=31, 0\"],\n  [16, 1, \"41\u00f79=4, 5\"],\n  [16, 2, \"30\u00f72=15, 0\"],\n  [16, 3, \"56\u00f76=9, 2\"],\n  [16, 4, \"99\u00f74=24, 3\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-fact answers in the single table on the page.\n# Every populated cell (5 per \"fact row\", at table rows 1, 5, 9, 13, 17 -\n# Word COM row/column indices are 1-based) gets its text replaced with a\n# newly generated problem/answer string.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Map of (tableRow, column) -> new cell text, in document order.\n# (row, col) use Word's 1-based Cell(row, col) addressing.\n$updates = @(\n    @(1, 1, \"89\u00f75=17, 4\"),\n    @(1, 2, \"33\u00f74=8, 1\"),\n    @(1, 3, \"96\u00f78=12, 0\"),\n    @(1, 4, \"74\u00f73=24, 2\"),\n    @(1, 5, \"50\u00f78=6, 2\"),\n\n    @(5, 1, \"89\u00f74=22, 1\"),\n    @(5, 2, \"44\u00f74=11, 0\"),\n    @(5, 3, \"22\u00f73=7, 1\"),\n    @(5, 4, \"31\u00f76=5, 1\"),\n    @(5, 5, \"60\u00f79=6, 6\"),\n\n    @(9, 1, \"10\u00f73=3, 1\"),\n    @(9, 2, \"53\u00f78=6, 5\"),\n    @(9, 3, \"24\u00f79=2, 6\"),\n    @(9, 4, \"17\u00f73=5, 2\"),\n    @(9, 5, \"92\u00f74=23, 0\"),\n\n    @(13, 1, \"73\u00f76=12, 1\"),\n    @(13, 2, \"94\u00f76=15, 4\"),\n    @(13, 3, \"88\u00f78=11, 0\"),\n    @(13, 4, \"32\u00f73=10, 2\"),\n    @(13, 5, \"92\u00f78=11, 4\"),\n\n    @(17, 1, \"93\u00f73=31, 0\"),\n    @(17, 2, \"41\u00f79=4, 5\"),\n    @(17, 3, \"30\u00f72=15, 0\"),\n    @(17, 4, \"56\u00f76=9, 2\"),\n    @(17, 5, \"99\u00f74=24, 3\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
